$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Change the Username and Password values to a masked placeholder
$ws.Range("B2").Value = "***"
$ws.Range("C2").Value = "***"

$wb.Save()
